$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.378.34'
$ws.Cells.Item(2, 5).Value = '  -2.07%  '

$ws.Cells.Item(3, 4).Value = '3.077.36'
$ws.Cells.Item(3, 5).Value = '  -2.23%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$__origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '525.28'
$ws.Cells.Item(5, 4).Style = $__origStyle
$ws.Cells.Item(5, 5).Value = '  -1.71%  '

$__origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '136.38'
$ws.Cells.Item(6, 4).Style = $__origStyle
$ws.Cells.Item(6, 5).Value = '  -4.96%  '

$ws.Cells.Item(7, 5).Value = '  -0.05%  '

$ws.Cells.Item(8, 4).Value = '3.076.65'
$ws.Cells.Item(8, 5).Value = '  -2.22%  '

$__origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.469'
$ws.Cells.Item(9, 4).Style = $__origStyle
$ws.Cells.Item(9, 5).Value = '  +4.48%  '

$__origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.24'
$ws.Cells.Item(10, 4).Style = $__origStyle
$ws.Cells.Item(10, 5).Value = '  +1.03%  '

$ws.Cells.Item(11, 5).Value = '  -3.50%  '

$__origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.401'
$ws.Cells.Item(12, 4).Style = $__origStyle
$ws.Cells.Item(12, 5).Value = '  +1.83%  '

$ws.Cells.Item(13, 5).Value = '  +1.66%  '

$ws.Cells.Item(14, 4).Value = '3.605.43'
$ws.Cells.Item(14, 5).Value = '  -2.17%  '

$__origStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '25.16'
$ws.Cells.Item(15, 4).Style = $__origStyle
$ws.Cells.Item(15, 5).Value = '  -2.26%  '

$__origStyle = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000160'
$ws.Cells.Item(16, 4).Style = $__origStyle
$ws.Cells.Item(16, 5).Value = '  -4.51%  '

$ws.Cells.Item(17, 4).Value = '57.395.26'
$ws.Cells.Item(17, 5).Value = '  -2.07%  '

$ws.Cells.Item(18, 4).Value = '3.075.75'
$ws.Cells.Item(18, 5).Value = '  -2.23%  '

$__origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '5.85'
$ws.Cells.Item(19, 4).Style = $__origStyle
$ws.Cells.Item(19, 5).Value = '  -4.68%  '

$__origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.45'
$ws.Cells.Item(20, 4).Style = $__origStyle
$ws.Cells.Item(20, 5).Value = '  -3.52%  '

$__origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.84'
$ws.Cells.Item(21, 4).Style = $__origStyle
$ws.Cells.Item(21, 5).Value = '  -2.03%  '

$__origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '348.43'
$ws.Cells.Item(22, 4).Style = $__origStyle
$ws.Cells.Item(22, 5).Value = '  +1.49%  '

$ws.Cells.Item(23, 5).Value = '  +0.10%  '

$__origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '68.49'
$ws.Cells.Item(24, 4).Style = $__origStyle
$ws.Cells.Item(24, 5).Value = '  +0.90%  '

$__origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.498'
$ws.Cells.Item(25, 4).Style = $__origStyle
$ws.Cells.Item(25, 5).Value = '  -3.52%  '

$__origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.166'
$ws.Cells.Item(26, 4).Style = $__origStyle
$ws.Cells.Item(26, 5).Value = '  -2.68%  '

$__origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 4).Style = $__origStyle
$ws.Cells.Item(27, 5).Value = '  +0.02%  '

$ws.Cells.Item(28, 4).Value = '0.0₃0845'
$ws.Cells.Item(28, 5).Value = '  -10.02%  '

$ws.Cells.Item(29, 5).Value = '  +0.07%  '

$__origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.11'
$ws.Cells.Item(30, 4).Style = $__origStyle
$ws.Cells.Item(30, 5).Value = '  -5.61%  '

$__origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.85'
$ws.Cells.Item(31, 4).Style = $__origStyle
$ws.Cells.Item(31, 5).Value = '  -2.46%  '

$__origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '5.88'
$ws.Cells.Item(32, 4).Style = $__origStyle
$ws.Cells.Item(32, 5).Value = '  -9.37%  '

$__origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '20.99'
$ws.Cells.Item(33, 4).Style = $__origStyle
$ws.Cells.Item(33, 5).Value = '  -1.18%  '

$ws.Cells.Item(34, 2).Value = 'NEARProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$__origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.82'
$ws.Cells.Item(34, 4).Style = $__origStyle
$ws.Cells.Item(34, 5).Value = '  -0.07%  '

$ws.Cells.Item(35, 2).Value = 'Monero'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '158.55'
$ws.Cells.Item(35, 4).Style = $__origStyle
$ws.Cells.Item(35, 5).Value = '  +0.08%  '

$__origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.12'
$ws.Cells.Item(36, 4).Style = $__origStyle
$ws.Cells.Item(36, 5).Value = '  -6.92%  '

$__origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.98'
$ws.Cells.Item(37, 4).Style = $__origStyle
$ws.Cells.Item(37, 5).Value = '  -4.51%  '

$__origStyle = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '25.48'
$ws.Cells.Item(38, 4).Style = $__origStyle
$ws.Cells.Item(38, 5).Value = '  -2.91%  '

$__origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.22'
$ws.Cells.Item(39, 4).Style = $__origStyle
$ws.Cells.Item(39, 5).Value = '  -4.55%  '

$ws.Cells.Item(40, 5).Value = '  -2.30%  '

$__origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.58'
$ws.Cells.Item(41, 4).Style = $__origStyle
$ws.Cells.Item(41, 5).Value = '  -5.39%  '

$__origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '4.02'
$ws.Cells.Item(42, 4).Style = $__origStyle
$ws.Cells.Item(42, 5).Value = '  -0.35%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '2.419.32'
$ws.Cells.Item(43, 5).Value = '  +4.29%  '

$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$__origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.688'
$ws.Cells.Item(44, 4).Style = $__origStyle
$ws.Cells.Item(44, 5).Value = '  -3.24%  '

$__origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '36.68'
$ws.Cells.Item(45, 4).Style = $__origStyle
$ws.Cells.Item(45, 5).Value = '  -0.32%  '

$ws.Cells.Item(46, 5).Value = '  +0.01%  '

$ws.Cells.Item(47, 4).Value = '3.115.29'
$ws.Cells.Item(47, 5).Value = '  -2.16%  '

$__origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0259'
$ws.Cells.Item(48, 4).Style = $__origStyle
$ws.Cells.Item(48, 5).Value = '  -2.80%  '

$__origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '5.97'
$ws.Cells.Item(49, 4).Style = $__origStyle
$ws.Cells.Item(49, 5).Value = '  -1.80%  '

$__origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.942'
$ws.Cells.Item(50, 4).Style = $__origStyle
$ws.Cells.Item(50, 5).Value = '  -7.16%  '

$__origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '19.35'
$ws.Cells.Item(51, 4).Style = $__origStyle
$ws.Cells.Item(51, 5).Value = '  -6.65%  '
